$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "36.612.96"
$ws.Cells.Item(2, 5).Value = "  +0.47%  "
$ws.Cells.Item(3, 4).Value = "1.964.26"
$ws.Cells.Item(3, 5).Value = "  +1.01%  "
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "'244.41"
$ws.Cells.Item(5, 5).Value = "  +0.50%  "
$ws.Cells.Item(6, 5).Value = "  +0.29%  "
$ws.Cells.Item(7, 4).Value = "'59.34"
$ws.Cells.Item(7, 5).Value = "  +1.75%  "
$ws.Cells.Item(8, 5).Value = "  -0.06%  "
$ws.Cells.Item(9, 4).Value = "'0.376"
$ws.Cells.Item(9, 5).Value = "  +2.86%  "
$ws.Cells.Item(10, 4).Value = "'0.0813"
$ws.Cells.Item(10, 5).Value = "  -2.39%  "
$ws.Cells.Item(11, 5).Value = "  -0.26%  "
$ws.Cells.Item(12, 4).Value = "'22.43"
$ws.Cells.Item(12, 5).Value = "  +3.66%  "
$ws.Cells.Item(13, 4).Value = "2.252.66"
$ws.Cells.Item(14, 4).Value = "'0.830"
$ws.Cells.Item(14, 5).Value = "  +1.15%  "
$ws.Cells.Item(15, 4).Value = "'13.77"
$ws.Cells.Item(15, 5).Value = "  +1.17%  "
$ws.Cells.Item(16, 5).Value = "  +0.61%  "
$ws.Cells.Item(17, 4).Value = "1.964.30"
$ws.Cells.Item(17, 5).Value = "  +1.07%  "
$ws.Cells.Item(18, 4).Value = "36.533.82"
$ws.Cells.Item(18, 5).Value = "  +0.42%  "
$ws.Cells.Item(19, 4).Value = "'69.99"
$ws.Cells.Item(19, 5).Value = "  +0.42%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0860"
$ws.Cells.Item(20, 5).Value = "  -0.21%  "
$ws.Cells.Item(21, 4).Value = "'229.33"
$ws.Cells.Item(21, 5).Value = "  -0.07%  "
$ws.Cells.Item(22, 4).Value = "'5.07"
$ws.Cells.Item(22, 5).Value = "  +0.48%  "
$ws.Cells.Item(23, 5).Value = "  -0.15%  "
$ws.Cells.Item(24, 5).Value = "  +0.66%  "
$ws.Cells.Item(25, 4).Value = "'2.36"
$ws.Cells.Item(25, 5).Value = "  +3.03%  "
$ws.Cells.Item(26, 4).Value = "'0.142"
$ws.Cells.Item(26, 5).Value = "  +7.71%  "
$ws.Cells.Item(27, 4).Value = "'9.23"
$ws.Cells.Item(27, 5).Value = "  +0.19%  "
$ws.Cells.Item(28, 4).Value = "'160.39"
$ws.Cells.Item(28, 5).Value = "  -1.11%  "
$ws.Cells.Item(29, 4).Value = "'19.45"
$ws.Cells.Item(29, 5).Value = "  +0.01%  "
$ws.Cells.Item(30, 4).Value = "'0.120"
$ws.Cells.Item(30, 5).Value = "  +1.84%  "
$ws.Cells.Item(31, 4).Value = "'1.16"
$ws.Cells.Item(31, 5).Value = "  +0.65%  "
$ws.Cells.Item(32, 4).Value = "'4.73"
$ws.Cells.Item(32, 5).Value = "  +1.40%  "
$ws.Cells.Item(33, 4).Value = "'0.0620"
$ws.Cells.Item(33, 5).Value = "  -1.26%  "
$ws.Cells.Item(34, 4).Value = "'4.29"
$ws.Cells.Item(34, 5).Value = "  +0.58%  "
$ws.Cells.Item(35, 5).Value = "  -0.02%  "
$ws.Cells.Item(36, 5).Value = "  +6.13%  "
$ws.Cells.Item(37, 4).Value = "'5.98"
$ws.Cells.Item(37, 5).Value = "  -4.40%  "
$ws.Cells.Item(38, 4).Value = "'3.37"
$ws.Cells.Item(38, 5).Value = "  +11.71%  "
$ws.Cells.Item(39, 5).Value = "  +0.21%  "
$ws.Cells.Item(40, 4).Value = "'0.0986"
$ws.Cells.Item(40, 5).Value = "  +0.86%  "
$ws.Cells.Item(41, 5).Value = "  +1.42%  "
$ws.Cells.Item(42, 5).Value = "  +0.20%  "
$ws.Cells.Item(43, 5).Value = "  +1.22%  "
$ws.Cells.Item(44, 4).Value = "'16.14"
$ws.Cells.Item(44, 5).Value = "  +0.86%  "
$ws.Cells.Item(45, 4).Value = "1.362.34"
$ws.Cells.Item(46, 5).Value = "  +0.82%  "
$ws.Cells.Item(47, 4).Value = "'87.91"
$ws.Cells.Item(47, 5).Value = "  +0.19%  "
$ws.Cells.Item(48, 4).Value = "'7.15"
$ws.Cells.Item(48, 5).Value = "  +0.52%  "
$ws.Cells.Item(49, 5).Value = "  +0.98%  "
$ws.Cells.Item(50, 4).Value = "2.143.06"
$ws.Cells.Item(51, 4).Value = "'43.90"
$ws.Cells.Item(51, 5).Value = "  -3.35%  "
